$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.448.90'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '1.570.23'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '208.97'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '22.18'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '0.0592'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').Value = '0.0866'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '1.794.06'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '1.582.49'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '3.82'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '0.518'
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').Value = '63.75'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').Value = '27.457.60'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '213.68'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = '0.0₃0692'
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').Value = '152.46'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').Value = '15.01'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '1.382.59'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +1.51%  '
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('D36').Value = '0.952'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = '0.544'
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').Value = '0.827'
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '0.988'
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('D44').Value = '64.27'
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = '1.706.20'
$ws.Range('D48').Value = '85.66'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('D49').Value = '0.0₇0999'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').Value = '0.0496'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '0.0956'
$ws.Range('E51').Value = '  -1.73%  '
